$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C values (new "Contract" figures read from external source)
$values = @(45, 50, 55, 60, 42, 47, 48, 49, 51, 52, 44, 40, 46, 50, 56, 38, 37, 35, 29, 28)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

# Update selection to C21
$ws.Range("C21").Select()

# Update workbook view window size/position
$win = $excel.ActiveWindow
$win.Left = 3480
$win.Top = 5100
$win.Width = 21600
$win.Height = 11295
